$wb = $excel.ActiveWorkbook

# Add the new "Replace Substrings" worksheet after the last existing sheet
# ("Trim Whitespace"), matching the GOMS-style task sheets already present.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$newSheet.Name = "Replace Substrings"

# Match the font sizing used by the other task sheets (13pt body, bold 13pt
# for the header/footer rows) before toggling Bold so no stray intermediate
# font gets recorded in the style table.
$newSheet.Range("A1:C5").Font.Size = 13
$newSheet.Range("A6:B6").Font.Size = 13

$newSheet.Range("A1:C1").Font.Bold = $true
$newSheet.Range("A6:B6").Font.Bold = $true

# Header row
$newSheet.Range("A1").Value = "Action"
$newSheet.Range("B1").Value = "Time"
$newSheet.Range("C1").Value = "Content"

# Step rows
$newSheet.Range("A2").Value = "Upload CSV"
$newSheet.Range("B2").Value = "5 min"
$newSheet.Range("C2").Value = "df = pd.read_csv('file.csv')"

$newSheet.Range("A3").Value = "Identify Substrings"
$newSheet.Range("B3").Value = "3 min"
$newSheet.Range("C3").Value = "Use df['column'].unique() to find unique values"

$newSheet.Range("A4").Value = "Replace Substrings"
$newSheet.Range("B4").Value = "2 min"
$newSheet.Range("C4").Value = "df['column'].str.replace('old_substring', 'new_substring', regex=True)"

$newSheet.Range("A5").Value = "Verify Changes"
$newSheet.Range("B5").Value = "1 min"
$newSheet.Range("C5").Value = "df['column'].unique() to check replacements"

# Footer / total row
$newSheet.Range("A6").Value = "Overall"
$newSheet.Range("B6").Value = "11 min"

# Leave the cursor where the author left it when they last saved the sheet.
[void]$newSheet.Range("M17").Select()
